$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear A12 and B12 (previously held "ActivePowerSummary" / "Мощность")
$ws.Range("A12:B12").ClearContents()

# Update the active selection to F7
$ws.Range("F7").Select()
